# "actualizacion de nota standar" - fill in the standard/note numbering
# column (A) for the existing data rows (3-36) with a sequential index.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 3; $i -le 36; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 2
}

# Restore the view/selection state captured in the saved workbook
# (scrolled so row 22 is at the top, cell B15 active).
$excel.Goto($ws.Range("A22"), $true)
$ws.Range("B15").Select()
